$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 120.44444
$ws.Range("I33").Value = 120.44444
$ws.Range("K33").Value = 120.44444
$ws.Range("M33").Value = 108.55556
$ws.Range("H64").Value = 3375
$ws.Range("I64").Value = 3133.3333
$ws.Range("K64").Value = 3133.3333
$ws.Range("M64").Value = -2885.3333
$ws.Range("H67").Value = 3375
$ws.Range("I67").Value = 3133.3333
$ws.Range("K67").Value = 3133.3333
$ws.Range("M67").Value = -2275.3333
$ws.Range("H74").Value = 4566.4287
$ws.Range("I74").Value = 4394.2
$ws.Range("K74").Value = 4394.2
$ws.Range("M74").Value = -3458.2
$ws.Range("H77").Value = 4566.4287
$ws.Range("I77").Value = 4394.2
$ws.Range("K77").Value = 21971
$ws.Range("M77").Value = -17291
$ws.Range("H135").Value = 718.6842
$ws.Range("I135").Value = 691
$ws.Range("K135").Value = 6219
$ws.Range("M135").Value = -3684
$ws.Range("H137").Value = 48881.477
$ws.Range("I137").Value = 1017.9167
$ws.Range("K137").Value = 3053.7501
$ws.Range("M137").Value = -503.7501000000002
$ws.Range("H138").Value = 1585.73
$ws.Range("I138").Value = 1081.6285
$ws.Range("J138").Value = 1857.1692
$ws.Range("K138").Value = 3244.8855
$ws.Range("L138").Value = 5571.5076
$ws.Range("M138").Value = 1895.1145
$ws.Range("N138").Value = -15851.5076
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 1973.1666
$ws.Range("I63").Value = 1973.1666
$ws.Range("K63").Value = 1973.1666
$ws.Range("M63").Value = -1287.1666
$ws.Range("H66").Value = 1973.1666
$ws.Range("I66").Value = 1973.1666
$ws.Range("K66").Value = 9865.833000000001
$ws.Range("M66").Value = -6433.833000000001
$ws.Range("H74").Value = 558.4737
$ws.Range("I74").Value = 558.4737
$ws.Range("K74").Value = 558.4737
$ws.Range("M74").Value = 315.5263
$ws.Range("H77").Value = 558.4737
$ws.Range("I77").Value = 558.4737
$ws.Range("K77").Value = 2792.3685
$ws.Range("M77").Value = 1575.6315
$ws.Range("H88").Value = 2472.5789
$ws.Range("I88").Value = 1967.4
$ws.Range("J88").Value = 3033.889
$ws.Range("K88").Value = 1967.4
$ws.Range("L88").Value = 3033.889
$ws.Range("M88").Value = -1561.4
$ws.Range("N88").Value = -3845.889
$ws.Range("H91").Value = 2472.5789
$ws.Range("I91").Value = 1967.4
$ws.Range("J91").Value = 3033.889
$ws.Range("K91").Value = 1967.4
$ws.Range("L91").Value = 3033.889
$ws.Range("M91").Value = -563.4000000000001
$ws.Range("N91").Value = -5841.889
$ws.Range("H110").Value = 584.6429000000001
$ws.Range("I110").Value = 584.6429000000001
$ws.Range("K110").Value = 584.6429000000001
$ws.Range("M110").Value = 1460.3571
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1162.4814
$ws.Range("I99").Value = 1038.2174
$ws.Range("K99").Value = 1038.2174
$ws.Range("M99").Value = 459.7826
$ws.Range("H134").Value = 3636.9243
$ws.Range("I134").Value = 3621.647
$ws.Range("J134").Value = 3688.8667
$ws.Range("K134").Value = 10864.941
$ws.Range("L134").Value = 11066.6001
$ws.Range("M134").Value = -8329.940999999999
$ws.Range("N134").Value = -16136.6001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2187.125
$ws.Range("J31").Value = 3079.8
$ws.Range("L31").Value = 3079.8
$ws.Range("N31").Value = -3669.8
$ws.Range("H34").Value = 2187.125
$ws.Range("J34").Value = 3079.8
$ws.Range("L34").Value = 3079.8
$ws.Range("N34").Value = -3483.8
$ws.Range("H62").Value = 2828
$ws.Range("I62").Value = 2900
$ws.Range("K62").Value = 2900
$ws.Range("M62").Value = -2276
$ws.Range("H65").Value = 2828
$ws.Range("I65").Value = 2900
$ws.Range("K65").Value = 14500
$ws.Range("M65").Value = -11380
$ws.Range("H141").Value = 51543.39
$ws.Range("J141").Value = 50340.06
$ws.Range("L141").Value = 50340.06
$ws.Range("N141").Value = -60700.06
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 74363.87
$ws.Range("J113").Value = 1245.8182
$ws.Range("L113").Value = 3737.4546
$ws.Range("N113").Value = -8077.4546
$ws.Range("H131").Value = 23117.871
$ws.Range("J131").Value = 31058.957
$ws.Range("L131").Value = 93176.871
$ws.Range("N131").Value = -103256.871
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H53").Value = 25000
$ws.Range("J53").Value = 25000
$ws.Range("L53").Value = 25000
$ws.Range("N53").Value = -26262
$ws.Range("H70").Value = 4440.7144
$ws.Range("I70").Value = 4269.9165
$ws.Range("K70").Value = 4269.9165
$ws.Range("M70").Value = -3999.9165
$ws.Range("H73").Value = 4440.7144
$ws.Range("I73").Value = 4269.9165
$ws.Range("K73").Value = 4269.9165
$ws.Range("M73").Value = -3333.9165
$ws.Range("H80").Value = 3456.2856
$ws.Range("I80").Value = 3665.6667
$ws.Range("J80").Value = 2200
$ws.Range("K80").Value = 3665.6667
$ws.Range("L80").Value = 2200
$ws.Range("M80").Value = -2667.6667
$ws.Range("N80").Value = -4196
$ws.Range("H83").Value = 3456.2856
$ws.Range("I83").Value = 3665.6667
$ws.Range("J83").Value = 2200
$ws.Range("K83").Value = 18328.3335
$ws.Range("L83").Value = 11000
$ws.Range("M83").Value = -13336.3335
$ws.Range("N83").Value = -20984
$ws.Range("H132").Value = 1328344.5
$ws.Range("I132").Value = 1481054
$ws.Range("K132").Value = 4443162
$ws.Range("M132").Value = -4440632
$ws.Range("H135").Value = 55000
$ws.Range("J135").Value = 55000
$ws.Range("L135").Value = 55000
$ws.Range("N135").Value = -65140
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2280.5
$ws.Range("I61").Value = 2128.55
$ws.Range("K61").Value = 2128.55
$ws.Range("M61").Value = -1926.55
$ws.Range("H100").Value = 1519.25
$ws.Range("I100").Value = 1519.25
$ws.Range("K100").Value = 1519.25
$ws.Range("M100").Value = -978.25
$ws.Range("H113").Value = 2280.5
$ws.Range("I113").Value = 2128.55
$ws.Range("K113").Value = 2128.55
$ws.Range("M113").Value = 41.44999999999982
$ws.Range("H132").Value = 4021.853
$ws.Range("I132").Value = 3357.8572
$ws.Range("K132").Value = 10073.5716
$ws.Range("M132").Value = -7543.571599999999
$ws.Range("H135").Value = 53385.668
$ws.Range("J135").Value = 53385.668
$ws.Range("L135").Value = 53385.668
$ws.Range("N135").Value = -63525.668
$ws.Range("H136").Value = 3089.4285
$ws.Range("I136").Value = 2398.2222
$ws.Range("K136").Value = 7194.6666
$ws.Range("M136").Value = -4644.6666
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 711.9
$ws.Range("I107").Value = 534.1875
$ws.Range("K107").Value = 1602.5625
$ws.Range("M107").Value = 317.4375
$ws.Range("H113").Value = 976.8461
$ws.Range("I113").Value = 799.25
$ws.Range("K113").Value = 2397.75
$ws.Range("M113").Value = -227.75
$ws.Range("H126").Value = 1228.3715
$ws.Range("I126").Value = 983.6667
$ws.Range("J126").Value = 2054.25
$ws.Range("K126").Value = 2951.0001
$ws.Range("L126").Value = 6162.75
$ws.Range("M126").Value = -481.0001000000002
$ws.Range("N126").Value = -11102.75
$ws.Range("H136").Value = 11575115
$ws.Range("I136").Value = 19841990
$ws.Range("J136").Value = 1489.1
$ws.Range("K136").Value = 59525970
$ws.Range("L136").Value = 4467.299999999999
$ws.Range("M136").Value = -59523420
$ws.Range("N136").Value = -9567.299999999999

Write-Output "done"